$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the existing most-recent meeting),
# shifting all existing data rows down by one.
$ws.Rows("2:2").Insert()

# Populate the new row 2 with the newest meeting (第630回).
$ws.Range("A2").Value = "第630回"
$ws.Range("B2").Value = "2025年11月26日（令和7年11月26日）"
$ws.Range("C2").Value = "１調査実施小委員会からの報告について`n２入院について（その７）`n３個別事項について（その９）データ提出加算`n`n"
$ws.Range("D2").Value = "－"
$ws.Range("E2").Value = "資料`n`n"
$ws.Range("F2").Value = "－"

# The multi-line text in C2 triggers Excel's auto row-height; reset it so the
# new row matches the sheet's default (unmarked) row height like the rest.
$ws.Rows("2:2").AutoFit()
